$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-24 01:25:21"

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
